# Weekly data refresh: a new week's price record is inserted at the top of
# the "Vega Modelo de Temuco - Bruselas (repollito)" data block (row 54),
# pushing all subsequent rows (old 54..159) down by one (new 55..160).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 54, shifting existing rows 54-159 down to 55-160.
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with this week's record.
$ws.Range("A54").Value = 10
$ws.Range("B54").Value = "Vega Modelo de Temuco"
$ws.Range("C54").Value = "La Araucanía"
$ws.Range("D54").Value = 45044
$ws.Range("E54").Value = 9
$ws.Range("F54").Value = 100112035
$ws.Range("G54").Value = "Bruselas (repollito)"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 40
$ws.Range("K54").Value = 28000
$ws.Range("L54").Value = 28000
$ws.Range("M54").Value = 28000
$ws.Range("N54").Value = "$/malla 15 kilos"
$ws.Range("O54").Value = "Región Metropolitana"
$ws.Range("P54").Value = 1867
$ws.Range("Q54").Value = 15
$ws.Range("R54").Value = "Hortaliza"
